# Actualización automática 2025-07-08 17:00:08
$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M7").Value = 697.36

# Sheet "VENTA MENSUAL"
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F7").Value = 2296.73
$ws2.Range("F22").Value = 19303.1

# Sheet "CUMPLIMIENTO MENSUAL"
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 15947.15
$ws3.Range("E16").Value = 28319.09
$ws3.Range("F16").Value = 0.3602553548708903
$ws3.Range("D19").Value = 19303.1
$ws3.Range("E19").Value = 46074.89762291768
$ws3.Range("F19").Value = 0.2952537658209566
